# Insert a new data row at row 101 (pushing existing rows 101:133 down to 102:134)
# and populate it with a new weekly price observation, consistent with the
# other rows in this "Albahaca" subset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 101 and below down by one row.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record.
$ws.Cells.Item(101, 1).Value = 8
$ws.Cells.Item(101, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(101, 3).Value = "Coquimbo"
$ws.Cells.Item(101, 4).Value = 44876
$ws.Cells.Item(101, 5).Value = 4
$ws.Cells.Item(101, 6).Value = 100112052
$ws.Cells.Item(101, 7).Value = "Albahaca"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 1000
$ws.Cells.Item(101, 11).Value = 4000
$ws.Cells.Item(101, 12).Value = 4500
$ws.Cells.Item(101, 13).Value = 4250
$ws.Cells.Item(101, 14).Value = "`$/paquete"
$ws.Cells.Item(101, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(101, 16).Value = 4250
$ws.Cells.Item(101, 17).Value = 1
$ws.Cells.Item(101, 18).Value = "Hortaliza"
